# Replace trailing asterisk footnote markers with superscript letter "a"
# in the header row and update the footnote text accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Follicular" + [char]0x1D43
$ws.Range("C1").Value = "FV-PTC" + [char]0x1D43
$ws.Range("D1").Value = "Papillary" + [char]0x1D43
$ws.Range("E1").Value = "Total" + [char]0x1D43

$ws.Range("A9").Value = [char]0x1D43 + " All values displayed as mean ± SD for ratio continuous variables or n (%) for dichotomous categorical variables. Percentages for the variant columns were calculated in respect to total patients within a variant (i.e., within column), and percentages for the total column was calculated in respect to the population total."
